$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.352819
$ws.Range("H2").Value = 10.058457
$ws.Range("I2").Value = 0.02224149976981271
$ws.Range("J2").Value = 0.02224149976981271
$ws.Range("M2").Value = 12.40685866666667
$ws.Range("N2").Value = 37.220576
$ws.Range("O2").Value = 0.1720325859617629
$ws.Range("P2").Value = 0.1720325859617629
$ws.Range("Q2").Value = 41.59795146791467
$ws.Range("R2").Value = 374.381563211232
$ws.Range("S2").Value = 0.003826262721068835
$ws.Range("T2").Value = 0.003826262721068835
$ws.Range("G3").Value = 3.352819
$ws.Range("H3").Value = 10.058457
$ws.Range("I3").Value = 0.02224149976981271
$ws.Range("J3").Value = 0.02224149976981271
$ws.Range("O3").Value = 0.6097142007069145
$ws.Range("P3").Value = 0.6097142007069145
$ws.Range("Q3").Value = 147.4305672295243
$ws.Range("R3").Value = 1326.875105065719
$ws.Range("S3").Value = 0.01356095825467438
$ws.Range("T3").Value = 0.01356095825467438
$ws.Range("G4").Value = 3.352819
$ws.Range("H4").Value = 10.058457
$ws.Range("I4").Value = 0.02224149976981271
$ws.Range("J4").Value = 0.02224149976981271
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2182532133313226
$ws.Range("P4").Value = 0.2182532133313226
$ws.Range("Q4").Value = 52.774226028845
$ws.Range("R4").Value = 474.968034259605
$ws.Range("S4").Value = 0.004854278794069496
$ws.Range("T4").Value = 0.004854278794069496
$ws.Range("I5").Value = 0.8292884613633072
$ws.Range("J5").Value = 0.8292884613633072
$ws.Range("M5").Value = 12.40685866666667
$ws.Range("N5").Value = 37.220576
$ws.Range("O5").Value = 0.1720325859617629
$ws.Range("P5").Value = 0.1720325859617629
$ws.Range("Q5").Value = 1551.006070890649
$ws.Range("R5").Value = 13959.05463801584
$ws.Range("S5").Value = 0.1426646385165813
$ws.Range("T5").Value = 0.1426646385165813
$ws.Range("I6").Value = 0.8292884613633072
$ws.Range("J6").Value = 0.8292884613633072
$ws.Range("O6").Value = 0.6097142007069145
$ws.Range("P6").Value = 0.6097142007069145
$ws.Range("Q6").Value = 5497.042444126572
$ws.Range("R6").Value = 49473.38199713915
$ws.Range("S6").Value = 0.5056289513755958
$ws.Range("T6").Value = 0.5056289513755958
$ws.Range("I7").Value = 0.8292884613633072
$ws.Range("J7").Value = 0.8292884613633072
$ws.Range("N7").Value = 47.220765
$ws.Range("O7").Value = 0.2182532133313226
$ws.Range("P7").Value = 0.2182532133313226
$ws.Range("S7").Value = 0.1809948714711302
$ws.Range("T7").Value = 0.1809948714711302
$ws.Range("I8").Value = 0.1484700388668802
$ws.Range("J8").Value = 0.1484700388668802
$ws.Range("M8").Value = 12.40685866666667
$ws.Range("N8").Value = 37.220576
$ws.Range("O8").Value = 0.1720325859617629
$ws.Range("P8").Value = 0.1720325859617629
$ws.Range("Q8").Value = 277.6813405185173
$ws.Range("R8").Value = 2499.132064666656
$ws.Range("S8").Value = 0.02554168472411285
$ws.Range("T8").Value = 0.02554168472411285
$ws.Range("I9").Value = 0.1484700388668802
$ws.Range("J9").Value = 0.1484700388668802
$ws.Range("O9").Value = 0.6097142007069145
$ws.Range("P9").Value = 0.6097142007069145
$ws.Range("S9").Value = 0.09052429107664438
$ws.Range("T9").Value = 0.09052429107664438
$ws.Range("I10").Value = 0.1484700388668802
$ws.Range("J10").Value = 0.1484700388668802
$ws.Range("N10").Value = 47.220765
$ws.Range("O10").Value = 0.2182532133313226
$ws.Range("P10").Value = 0.2182532133313226
$ws.Range("Q10").Value = 352.2870072056349
$ws.Range("R10").Value = 3170.583064850714
$ws.Range("S10").Value = 0.03240406306612296
$ws.Range("T10").Value = 0.03240406306612296
